$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new "time_taken" column, styled like the other headers
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Data values for the new column
$times = @(
    "2021-10-05 13:39:54.814419",
    "2021-10-05 13:39:54.814430",
    "2021-10-05 13:39:54.814433",
    "2021-10-05 13:39:54.814436",
    "2021-10-05 13:39:54.814438",
    "2021-10-05 13:39:54.814441",
    "2021-10-05 13:39:54.814444",
    "2021-10-05 13:39:54.814446",
    "2021-10-05 13:39:54.814449",
    "2021-10-05 13:39:54.814451",
    "2021-10-05 13:39:54.814454",
    "2021-10-05 13:39:54.814456"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
